# "Updated summary report formatting"
#
# The original cell text was missing a space after the "$>" shell prompt;
# fix the wording and put the selection/active-cell back on A1 (the sheet's
# natural default) instead of the stray A2 that had been left selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the shared-string text: "$>echo ..." -> "$> echo ..."
$ws.Range("A1").Value = "$> echo 'Hello World' | grep Hello && cat < input.txt >> output.txt"

# Restore the default selection to A1 (was left on A2 in the saved file).
$ws.Range("A1").Select() | Out-Null
